# Questionnaire results - add percentage breakdown table (columns H:N)
# mirroring the raw counts table (columns A:F) on the worksheet named
# "Sheet1" (tabSelected / active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Pad the raw-count table so every question row has a value in B:F ---
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# --- Mirror headers / question labels into H:M, H2:H6 (reuse shared strings) ---
$ws.Range("H1").Value = $ws.Range("A1").Value()
$ws.Range("I1").Value = $ws.Range("B1").Value()
$ws.Range("J1").Value = $ws.Range("C1").Value()
$ws.Range("K1").Value = $ws.Range("D1").Value()
$ws.Range("L1").Value = $ws.Range("E1").Value()
$ws.Range("M1").Value = $ws.Range("F1").Value()

$ws.Range("H2").Value = $ws.Range("A2").Value()
$ws.Range("H3").Value = $ws.Range("A3").Value()
$ws.Range("H4").Value = $ws.Range("A4").Value()
$ws.Range("H5").Value = $ws.Range("A5").Value()
$ws.Range("H6").Value = $ws.Range("A6").Value()

# --- Row totals (column N) ---
$ws.Range("N2").Formula = "=SUM(B2:F2)"

# --- Row 2 percentages (not part of the shared-formula group below) ---
$ws.Range("I2").Formula = "=B2/N2"
$ws.Range("J2").Formula = "=C2/N2"
$ws.Range("K2").Formula = "=D2/N2"
$ws.Range("L2").Formula = "=E2/N2"
$ws.Range("M2").Formula = "=F2/N2"

# --- Rows 3-6 percentages + totals, filled as shared formulas (si order: I,J,K,L,M,N) ---
$ws.Range("I3:I6").Formula = "=B3/N3"
$ws.Range("J3:J6").Formula = "=C3/N3"
$ws.Range("K3:K6").Formula = "=D3/N3"
$ws.Range("L3:L6").Formula = "=E3/N3"
$ws.Range("M3:M6").Formula = "=F3/N3"
$ws.Range("N3:N6").Formula = "=SUM(B3:F3)"

# --- Format the percentage columns as percentages with 2 decimals (numFmtId 10) ---
$ws.Range("I2:M6").NumberFormat = "0.00%"

# --- Column widths for the mirrored columns (match A and F) ---
$ws.Columns.Item(8).ColumnWidth = 45.16666666666667
$ws.Columns.Item(13).ColumnWidth = 15.736979166666666

# --- Selection / active cell moved to L15 ---
$ws.Range("L15").Select()

# --- Page setup: paper size 9 (A4), portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
